$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 300
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 300
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 300
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -638
$ws.Range("H51").Value = 3159
$ws.Range("J51").Value = 3448.75
$ws.Range("L51").Value = 3448.75
$ws.Range("N51").Value = -4416.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 3340000
$ws.Range("J6").Value = 10000
$ws.Range("L6").Value = 10000
$ws.Range("N6").Value = -10346
$ws.Range("H26").Value = 5002.3335
$ws.Range("I26").Value = 5002.3335
$ws.Range("K26").Value = 5002.3335
$ws.Range("M26").Value = -4672.3335
$ws.Range("H32").Value = 5052.16
$ws.Range("I32").Value = 4117.477
$ws.Range("K32").Value = 4117.477
$ws.Range("M32").Value = -3830.477
$ws.Range("H45").Value = 1365.619
$ws.Range("I45").Value = 944.1818
$ws.Range("K45").Value = 944.1818
$ws.Range("M45").Value = -567.1818
$ws.Range("H61").Value = 5532.846
$ws.Range("I61").Value = 7026.1333
$ws.Range("K61").Value = 7026.1333
$ws.Range("M61").Value = -6814.1333
$ws.Range("H74").Value = 1755.2084
$ws.Range("I74").Value = 485
$ws.Range("K74").Value = 485
$ws.Range("M74").Value = 389
$ws.Range("H77").Value = 1755.2084
$ws.Range("I77").Value = 485
$ws.Range("K77").Value = 2425
$ws.Range("M77").Value = 1943
$ws.Range("H122").Value = 1603.1538
$ws.Range("I122").Value = 1171.8572
$ws.Range("K122").Value = 3515.5716
$ws.Range("M122").Value = -1065.5716
$ws.Range("H132").Value = 2278.25
$ws.Range("I132").Value = 1915.7646
$ws.Range("J132").Value = 4332.3335
$ws.Range("K132").Value = 5747.293799999999
$ws.Range("L132").Value = 12997.0005
$ws.Range("M132").Value = -3217.293799999999
$ws.Range("N132").Value = -18057.0005
$ws.Range("H136").Value = 5532.846
$ws.Range("I136").Value = 7026.1333
$ws.Range("K136").Value = 21078.3999
$ws.Range("M136").Value = -18528.3999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 69000
$ws.Range("J13").Value = 69000
$ws.Range("L13").Value = 69000
$ws.Range("N13").Value = -69336
$ws.Range("H20").Value = 3235.2
$ws.Range("I20").Value = 2829.3333
$ws.Range("K20").Value = 2829.3333
$ws.Range("M20").Value = -2582.3333
$ws.Range("H105").Value = 2634.6191
$ws.Range("I105").Value = 2333
$ws.Range("K105").Value = 2333
$ws.Range("M105").Value = -586
$ws.Range("H107").Value = 974.6111
$ws.Range("I107").Value = 808.2222
$ws.Range("J107").Value = 1141
$ws.Range("K107").Value = 808.2222
$ws.Range("L107").Value = 1141
$ws.Range("M107").Value = 1111.7778
$ws.Range("N107").Value = -4981
$ws.Range("H134").Value = 1275.1538
$ws.Range("I134").Value = 999.85297
$ws.Range("K134").Value = 2999.55891
$ws.Range("M134").Value = -464.5589100000002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 291.875
$ws.Range("I2").Value = 283.75
$ws.Range("K2").Value = 1702.5
$ws.Range("M2").Value = -1589.5
$ws.Range("H16").Value = 995
$ws.Range("I16").Value = 995
$ws.Range("K16").Value = 2985
$ws.Range("M16").Value = -2812
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("H21").Value = 900
$ws.Range("J21").Value = 900
$ws.Range("L21").Value = 2700
$ws.Range("N21").Value = -3046
$ws.Range("H22").Value = 3232.6667
$ws.Range("I22").Value = 3000
$ws.Range("J22").Value = 3253.818
$ws.Range("K22").Value = 9000
$ws.Range("L22").Value = 9761.454000000002
$ws.Range("M22").Value = -8831
$ws.Range("N22").Value = -10099.454
$ws.Range("H26").Value = 496.66666
$ws.Range("I26").Value = 500
$ws.Range("J26").Value = 495
$ws.Range("K26").Value = 1500
$ws.Range("L26").Value = 1485
$ws.Range("M26").Value = -1212
$ws.Range("N26").Value = -2061
$ws.Range("H27").Value = 3232.6667
$ws.Range("I27").Value = 3000
$ws.Range("J27").Value = 3253.818
$ws.Range("K27").Value = 9000
$ws.Range("L27").Value = 9761.454000000002
$ws.Range("M27").Value = -8898
$ws.Range("N27").Value = -9965.454000000002
$ws.Range("H33").Value = 122.083336
$ws.Range("I33").Value = 114.1
$ws.Range("K33").Value = 684.5999999999999
$ws.Range("M33").Value = -401.5999999999999
$ws.Range("H34").Value = 8846
$ws.Range("I34").Value = 12616.25
$ws.Range("J34").Value = 1305.5
$ws.Range("K34").Value = 37848.75
$ws.Range("L34").Value = 3916.5
$ws.Range("M34").Value = -37764.75
$ws.Range("N34").Value = -4084.5
$ws.Range("H39").Value = 1985.2
$ws.Range("J39").Value = 2306.5
$ws.Range("L39").Value = 6919.5
$ws.Range("N39").Value = -7507.5
$ws.Range("H40").Value = 104.75
$ws.Range("J40").Value = 130
$ws.Range("L40").Value = 520
$ws.Range("N40").Value = -658
$ws.Range("H50").Value = 142958300
$ws.Range("J50").Value = 200002080
$ws.Range("L50").Value = 600006240
$ws.Range("N50").Value = -600007202
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()
$ws.Range("H53").Value = 142958300
$ws.Range("J53").Value = 200002080
$ws.Range("L53").Value = 600006240
$ws.Range("N53").Value = -600007202
$ws.Range("H57").Value = 2100
$ws.Range("J57").Value = 4000
$ws.Range("L57").Value = 12000
$ws.Range("N57").Value = -13118
$ws.Range("H58").Value = 2500
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H131").Value = 13985.905
$ws.Range("J131").Value = 15069.857
$ws.Range("L131").Value = 45209.571
$ws.Range("N131").Value = -55289.571

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2505850
$ws.Range("H46").Value = 23800
$ws.Range("J46").Value = 24222.223
$ws.Range("L46").Value = 24222.223
$ws.Range("N46").Value = -24534.223
$ws.Range("H80").Value = 2278.6
$ws.Range("I80").Value = 2490.1428
$ws.Range("K80").Value = 2490.1428
$ws.Range("M80").Value = -1492.1428
$ws.Range("H83").Value = 2278.6
$ws.Range("I83").Value = 2490.1428
$ws.Range("K83").Value = 12450.714
$ws.Range("M83").Value = -7458.714
$ws.Range("H102").Value = 1778.8387
$ws.Range("I102").Value = 1391.0416
$ws.Range("K102").Value = 1391.0416
$ws.Range("M102").Value = 230.9584

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H46").Value = 1658.125
$ws.Range("J46").Value = 2243.875
$ws.Range("L46").Value = 2243.875
$ws.Range("N46").Value = -2619.875
$ws.Range("H61").Value = 2099.2727
$ws.Range("I61").Value = 2022.9048
$ws.Range("K61").Value = 2022.9048
$ws.Range("M61").Value = -1820.9048
$ws.Range("H113").Value = 2099.2727
$ws.Range("I113").Value = 2022.9048
$ws.Range("K113").Value = 2022.9048
$ws.Range("M113").Value = 147.0952
$ws.Range("H136").Value = 1953.6666
$ws.Range("I136").Value = 1575.125
$ws.Range("K136").Value = 4725.375
$ws.Range("M136").Value = -2175.375

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 807.3333
$ws.Range("I107").Value = 598.3333
$ws.Range("K107").Value = 1794.9999
$ws.Range("M107").Value = 125.0001
$ws.Range("H122").Value = 34323.707
$ws.Range("I122").Value = 66854.25
$ws.Range("J122").Value = 1793.1666
$ws.Range("K122").Value = 200562.75
$ws.Range("L122").Value = 5379.4998
$ws.Range("M122").Value = -198112.75
$ws.Range("N122").Value = -10279.4998
$ws.Range("H126").Value = 4998.7407
$ws.Range("I126").Value = 7592.75
$ws.Range("J126").Value = 1225.6364
$ws.Range("K126").Value = 22778.25
$ws.Range("L126").Value = 3676.9092
$ws.Range("M126").Value = -20308.25
$ws.Range("N126").Value = -8616.9092
$ws.Range("H132").Value = 1944.8206
$ws.Range("I132").Value = 1232.5333
$ws.Range("K132").Value = 3697.5999
$ws.Range("M132").Value = -1167.5999
$ws.Range("H136").Value = 1240.7273
$ws.Range("I136").Value = 831.56525
$ws.Range("K136").Value = 2494.69575
$ws.Range("M136").Value = 55.30425000000014
